# Updates cryptos list: refreshed Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.334.23"
$ws.Range("E2").Value = "  -0.91%  "
# Row 3
$ws.Range("D3").Value = "3.222.26"
$ws.Range("E3").Value = "  -1.37%  "
# Row 4
$ws.Range("E4").Value = "  -0.02%  "
# Row 5
$ws.Range("D5").Value = "'578.13"
$ws.Range("E5").Value = "  -1.45%  "
# Row 6
$ws.Range("D6").Value = "'183.52"
$ws.Range("E6").Value = "  -1.56%  "
# Row 7
$ws.Range("E7").Value = "  +0.00%  "
# Row 8
$ws.Range("D8").Value = "'0.607"
$ws.Range("E8").Value = "  +1.15%  "
# Row 9
$ws.Range("D9").Value = "3.218.21"
$ws.Range("E9").Value = "  -1.53%  "
# Row 10
$ws.Range("E10").Value = "  -2.93%  "
# Row 11
$ws.Range("E11").Value = "  -2.41%  "
# Row 12
$ws.Range("E12").Value = "  -1.71%  "
# Row 13
$ws.Range("D13").Value = "3.776.36"
$ws.Range("E13").Value = "  -1.65%  "
# Row 14
$ws.Range("E14").Value = "  +0.14%  "
# Row 15
$ws.Range("E15").Value = "  -3.20%  "
# Row 16
$ws.Range("D16").Value = "67.394.31"
$ws.Range("E16").Value = "  -0.88%  "
# Row 17
$ws.Range("E17").Value = "  -2.04%  "
# Row 18
$ws.Range("D18").Value = "3.208.87"
$ws.Range("E18").Value = "  -1.88%  "
# Row 19
$ws.Range("D19").Value = "'5.74"
$ws.Range("E19").Value = "  -2.05%  "
# Row 20
$ws.Range("E20").Value = "  -1.60%  "
# Row 21
$ws.Range("D21").Value = "'394.62"
$ws.Range("E21").Value = "  +3.24%  "
# Row 22
$ws.Range("E22").Value = "  -2.41%  "
# Row 23
$ws.Range("E23").Value = "  +0.23%  "
# Row 24
$ws.Range("E24").Value = "  -0.40%  "
# Row 25
$ws.Range("E25").Value = "  -0.20%  "
# Row 26
$ws.Range("E26").Value = "  -2.94%  "
# Row 27
$ws.Range("E27").Value = "  -1.99%  "
# Row 28
$ws.Range("E28").Value = "  -3.25%  "
# Row 29
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.12%  "
# Row 30
$ws.Range("E30").Value = "  -2.41%  "
# Row 31
$ws.Range("E31").Value = "  -5.59%  "
# Row 32
$ws.Range("D32").Value = "'22.51"
$ws.Range("E32").Value = "  -1.46%  "
# Row 33
$ws.Range("E33").Value = "  -3.58%  "
# Row 35
$ws.Range("D35").Value = "'1.25"
$ws.Range("E35").Value = "  -2.43%  "
# Row 36
$ws.Range("D36").Value = "'160.29"
$ws.Range("E36").Value = "  -1.28%  "
# Row 37
$ws.Range("D37").Value = "'1.47"
$ws.Range("E37").Value = "  -4.66%  "
# Row 38
$ws.Range("E38").Value = "  +1.18%  "
# Row 39
$ws.Range("E39").Value = "  -0.87%  "
# Row 40
$ws.Range("E40").Value = "  -4.49%  "
# Row 41
$ws.Range("E41").Value = "  -1.50%  "
# Row 42
$ws.Range("D42").Value = "'6.50"
$ws.Range("E42").Value = "  -4.26%  "
# Row 43
$ws.Range("E43").Value = "  -6.15%  "
# Row 44
$ws.Range("D44").Value = "'0.0682"
$ws.Range("E44").Value = "  -1.36%  "
# Row 45
$ws.Range("D45").Value = "'40.40"
$ws.Range("E45").Value = "  -2.21%  "
# Row 46
$ws.Range("D46").Value = "2.590.80"
$ws.Range("E46").Value = "  -2.15%  "
# Row 47
$ws.Range("D47").Value = "'333.43"
$ws.Range("E47").Value = "  -2.60%  "
# Row 48
$ws.Range("D48").Value = "'24.45"
$ws.Range("E48").Value = "  -3.76%  "
# Row 49
$ws.Range("D49").Value = "'0.0277"
$ws.Range("E49").Value = "  -2.76%  "
# Row 50
$ws.Range("E50").Value = "  -0.15%  "
# Row 51
$ws.Range("E51").Value = "  -1.34%  "
